{"js": "// no-op test\nawait context.sync();\n", "ps1": "# no-op test\n$d = $word.ActiveDocument\n"}
